$d = $word.ActiveDocument

# 1. Find the paragraph starting with "Email.Username" and insert a _GoBack bookmark
#    right before it (at the very start of the paragraph, before the proofErr/run).
$rng = $d.Content
$found = $rng.Find.Execute("Email.Username", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPoint = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null
}

# 2. Remove the trailing "FilePattern_org = ..." paragraph entirely (including the
#    _GoBack bookmark that used to sit at the end of the document).
$paraRange = $d.Content
$found2 = $paraRange.Find.Execute("FilePattern_org", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Expand to the full paragraph containing this text
    $para = $paraRange.Paragraphs(1).Range
    $para.Delete()
}

$d.Save()
